# Scheduled runner update: refresh currentAveragePrice / Leve profit
# figures across several item sheets (ALC, ARM, CRP, CUL, GSM, LTW, WVR)
# to their latest market-board snapshot values.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 1807.4445
$ws.Range("I33").Value = 783.375
$ws.Range("J33").Value = 10000
$ws.Range("K33").Value = 783.375
$ws.Range("L33").Value = 10000
$ws.Range("M33").Value = -554.375
$ws.Range("N33").Value = -10458
$ws.Range("H55").Value = 59682.766
$ws.Range("I55").Value = 320
$ws.Range("J55").Value = 67597.8
$ws.Range("K55").Value = 320
$ws.Range("L55").Value = 67597.8
$ws.Range("M55").Value = -106
$ws.Range("N55").Value = -68025.8
$ws.Range("H127").Value = 17629.572
$ws.Range("I127").Value = 22681.4
$ws.Range("K127").Value = 68044.20000000001
$ws.Range("M127").Value = -63084.20000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3499.582
$ws.Range("I32").Value = 3133.9355
$ws.Range("K32").Value = 3133.9355
$ws.Range("M32").Value = -2846.9355
$ws.Range("H61").Value = 3772.15
$ws.Range("I61").Value = 2797.0344
$ws.Range("K61").Value = 2797.0344
$ws.Range("M61").Value = -2585.0344
$ws.Range("H97").Value = 1749.2084
$ws.Range("I97").Value = 2074.9443
$ws.Range("J97").Value = 772
$ws.Range("K97").Value = 2074.9443
$ws.Range("L97").Value = 772
$ws.Range("M97").Value = -1578.9443
$ws.Range("N97").Value = -1764
$ws.Range("H132").Value = 4472.4814
$ws.Range("I132").Value = 1849.1666
$ws.Range("K132").Value = 5547.4998
$ws.Range("M132").Value = -3017.4998
$ws.Range("H136").Value = 3772.15
$ws.Range("I136").Value = 2797.0344
$ws.Range("K136").Value = 8391.1032
$ws.Range("M136").Value = -5841.1032

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 241536.36
$ws.Range("I31").Value = 418329.16
$ws.Range("K31").Value = 418329.16
$ws.Range("M31").Value = -418034.16
$ws.Range("H34").Value = 241536.36
$ws.Range("I34").Value = 418329.16
$ws.Range("K34").Value = 418329.16
$ws.Range("M34").Value = -418127.16
$ws.Range("H94").Value = 754.9091
$ws.Range("I94").Value = 552.25
$ws.Range("K94").Value = 552.25
$ws.Range("M94").Value = -101.25
$ws.Range("H132").Value = 3544.0444
$ws.Range("I132").Value = 2729.7334
$ws.Range("K132").Value = 8189.2002
$ws.Range("M132").Value = -5659.2002
$ws.Range("H134").Value = 199821.6
$ws.Range("I134").Value = 2560.8286
$ws.Range("K134").Value = 7682.485799999999
$ws.Range("M134").Value = -5147.485799999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H41").Value = 334.25
$ws.Range("I41").Value = 334.25
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 1002.75
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = -664.75
$ws.Range("N41").ClearContents()
$ws.Range("H122").Value = 500502.5
$ws.Range("I122").Value = 1000
$ws.Range("J122").Value = 1000005
$ws.Range("K122").Value = 9000
$ws.Range("L122").Value = 9000045
$ws.Range("M122").Value = -6550
$ws.Range("N122").Value = -9004945
$ws.Range("H128").Value = 299997.5
$ws.Range("I128").Value = 299997.5
$ws.Range("K128").Value = 899992.5
$ws.Range("M128").Value = -895012.5
$ws.Range("H129").Value = 2277.077
$ws.Range("J129").Value = 2773.5789
$ws.Range("L129").Value = 8320.736699999999
$ws.Range("N129").Value = -18320.7367
$ws.Range("H132").Value = 3020.6428
$ws.Range("I132").Value = 1365.4445
$ws.Range("J132").Value = 6000
$ws.Range("K132").Value = 12289.0005
$ws.Range("L132").Value = 54000
$ws.Range("M132").Value = -9759.0005
$ws.Range("N132").Value = -59060

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2138.8235
$ws.Range("I102").Value = 1742
$ws.Range("K102").Value = 1742
$ws.Range("M102").Value = -120
$ws.Range("H122").Value = 4621.75
$ws.Range("I122").Value = 4013.2856
$ws.Range("J122").Value = 6041.5
$ws.Range("K122").Value = 12039.8568
$ws.Range("L122").Value = 18124.5
$ws.Range("M122").Value = -9589.856800000001
$ws.Range("N122").Value = -23024.5
$ws.Range("H132").Value = 526408.5600000001
$ws.Range("I132").Value = 528872.6
$ws.Range("J132").Value = 503000
$ws.Range("K132").Value = 1586617.8
$ws.Range("L132").Value = 1509000
$ws.Range("M132").Value = -1584087.8
$ws.Range("N132").Value = -1514060
$ws.Range("H135").Value = 99999.8
$ws.Range("J135").Value = 99999.8
$ws.Range("L135").Value = 99999.8
$ws.Range("N135").Value = -110139.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1814.16
$ws.Range("I22").Value = 1760.5264
$ws.Range("J22").Value = 1984
$ws.Range("K22").Value = 1760.5264
$ws.Range("L22").Value = 1984
$ws.Range("M22").Value = -1465.5264
$ws.Range("N22").Value = -2574
$ws.Range("H27").Value = 1814.16
$ws.Range("I27").Value = 1760.5264
$ws.Range("J27").Value = 1984
$ws.Range("K27").Value = 1760.5264
$ws.Range("L27").Value = 1984
$ws.Range("M27").Value = -1653.5264
$ws.Range("N27").Value = -2198
$ws.Range("H55").Value = 1025.5454
$ws.Range("I55").Value = 444.26923
$ws.Range("J55").Value = 3184.5715
$ws.Range("K55").Value = 444.26923
$ws.Range("L55").Value = 3184.5715
$ws.Range("M55").Value = -271.26923
$ws.Range("N55").Value = -3530.5715
$ws.Range("H68").Value = 4423.875
$ws.Range("I68").Value = 4448.857
$ws.Range("K68").Value = 4448.857
$ws.Range("M68").Value = -3699.857
$ws.Range("H71").Value = 4423.875
$ws.Range("I71").Value = 4448.857
$ws.Range("K71").Value = 22244.285
$ws.Range("M71").Value = -18500.285
$ws.Range("H122").Value = 1056897.9
$ws.Range("I122").Value = 837004.2
$ws.Range("K122").Value = 2511012.6
$ws.Range("M122").Value = -2508562.6
$ws.Range("H132").Value = 5145.4443
$ws.Range("J132").Value = 5264.4287
$ws.Range("L132").Value = 15793.2861
$ws.Range("N132").Value = -20853.2861
$ws.Range("H136").Value = 3258.8708
$ws.Range("I136").Value = 2473.9092
$ws.Range("J136").Value = 5177.6665
$ws.Range("K136").Value = 7421.7276
$ws.Range("L136").Value = 15532.9995
$ws.Range("M136").Value = -4871.7276
$ws.Range("N136").Value = -20632.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 25555
$ws.Range("J15").Value = 25555
$ws.Range("L15").Value = 25555
$ws.Range("N15").Value = -26131
$ws.Range("H100").Value = 454.35294
$ws.Range("I100").Value = 406.0909
$ws.Range("K100").Value = 812.1818
$ws.Range("M100").Value = -271.1818
